$d = $word.ActiveDocument
$rng = $d.Content
$rng.Find.Execute("<Proposal Description>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para = $rng.Paragraphs(1)
$startRng = $d.Range($para.Range.Start, $para.Range.Start)
$startRng.InsertParagraphBefore()

$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$rng2 = $d.Content
$rng2.Find.Execute("<Proposal Description>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $rng2.Paragraphs(1)
$newStart = $d.Range($para2.Range.Start, $para2.Range.Start)
$d.Bookmarks.Add("_GoBack", $newStart)

$d.Repaginate()
